$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 516, shifting existing rows 516..571 down to 517..572
$ws.Rows.Item(516).Insert()

# Populate the newly inserted row 516 with the new data record
$ws.Cells.Item(516, 1).Value = 5
$ws.Cells.Item(516, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(516, 3).Value = "Maule"
$ws.Cells.Item(516, 4).Value = 44946
$ws.Cells.Item(516, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(516, 5).Value = 7
$ws.Cells.Item(516, 6).Value = 100112043
$ws.Cells.Item(516, 7).Value = "Pepino ensalada"
$ws.Cells.Item(516, 8).Value = "Sin especificar"
$ws.Cells.Item(516, 9).Value = "Primera"
$ws.Cells.Item(516, 10).Value = 700
$ws.Cells.Item(516, 11).Value = 7000
$ws.Cells.Item(516, 12).Value = 8000
$ws.Cells.Item(516, 13).Value = 7286
$ws.Cells.Item(516, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(516, 15).Value = "Región del Maule"
$ws.Cells.Item(516, 16).Value = 91
$ws.Cells.Item(516, 17).Value = 80
$ws.Cells.Item(516, 18).Value = "Hortaliza"
